# 04/01/26 New Year Push
# Add a new "TestcaseID" column at the front of the test data sheet and
# clear out the now-redundant RESPONSEPAYLOAD values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column before column A, shifting the existing
# REQUESTPAYLOAD / key / address / RESPONSEPAYLOAD / place_id / STATUS
# columns one place to the right.
$ws.Columns.Item(1).Insert()
$ws.Columns.Item(1).ColumnWidth = 19.75

# Populate the new TestcaseID column.
$ws.Cells.Item(1, 1).Value = "TestcaseID"
$ws.Cells.Item(2, 1).Value = "TS01_TC01"
$ws.Cells.Item(3, 1).Value = "TS02_TC01"

# The RESPONSEPAYLOAD column (now column E) used to just duplicate the
# address column values with a highlighted style - clear it out and
# drop back to the normal style.
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 5).Style = "Normal"

# The place_id column (now column F) was also using the highlighted
# style - restore it back to the normal style.
$ws.Cells.Item(2, 6).Style = "Normal"
$ws.Cells.Item(3, 6).Style = "Normal"

# Move the active selection to A3, matching the saved selection state.
$ws.Range("A3").Select() | Out-Null
